# Update uncertainty-related constants on the "constants" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

$ws.Range("B2").Value  = 10.19985961024888
$ws.Range("B4").Value  = 1856.545344897319
$ws.Range("B5").Value  = 24358.99779796487
$ws.Range("B8").Value  = 0.809744300858402
$ws.Range("B9").Value  = 0.7073697945015182
$ws.Range("B10").Value = 2.822627157457899
$ws.Range("B11").Value = 0.3840790815586094
